# Update "想去人数" (number of people wanting to go) counts for several
# rows across three worksheets, per the regenerated site data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 863
$ws1.Range("F6").Value  = 670
$ws1.Range("F7").Value  = 1240
$ws1.Range("F9").Value  = 830
$ws1.Range("F15").Value = 964
$ws1.Range("F16").Value = 10073
$ws1.Range("F17").Value = 635
$ws1.Range("F25").Value = 291
$ws1.Range("F29").Value = 280
$ws1.Range("F30").Value = 193
$ws1.Range("F32").Value = 72
$ws1.Range("F35").Value = 180
$ws1.Range("F36").Value = 199

# --- Sheet "本地生活" (local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 823

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 823
$ws4.Range("F8").Value  = 863
$ws4.Range("F9").Value  = 670
$ws4.Range("F10").Value = 1240
$ws4.Range("F14").Value = 830
$ws4.Range("F19").Value = 964
$ws4.Range("F20").Value = 10073
$ws4.Range("F22").Value = 635
$ws4.Range("F36").Value = 280
$ws4.Range("F37").Value = 193
$ws4.Range("F39").Value = 72
$ws4.Range("F43").Value = 180
$ws4.Range("F46").Value = 199
